# NATMI LR-pair table (Sema3f-Nrp2): add the 'M2' sending-cluster rows and
# refresh the per-cluster expression/specificity statistics for every
# Sending x Target cluster combination (now 4 x 4 = 16 data rows instead of
# 3 x 4 = 12), per Dr Hou's advice.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A Sending cluster | B Ligand symbol | C Receptor symbol |
# D Target cluster | E..T the numeric NATMI statistics.
$arr = New-Object 'object[,]' 16,20
# row 2: ECs -> ECs
$arr[0,0] = 'ECs'
$arr[0,1] = 'Sema3f'
$arr[0,2] = 'Nrp2'
$arr[0,3] = 'ECs'
$arr[0,4] = 3
$arr[0,5] = 1
$arr[0,6] = 37.05582933333334
$arr[0,7] = 111.167488
$arr[0,8] = 0.9062076988734117
$arr[0,9] = 0.9062076988734117
$arr[0,10] = 3
$arr[0,11] = 1
$arr[0,12] = 31.632955
$arr[0,13] = 94.898865
$arr[0,14] = 0.5000340016649593
$arr[0,15] = 0.5000340016649593
$arr[0,16] = 1172.185381789013
$arr[0,17] = 10549.66843610112
$arr[0,18] = 0.4531346620072664
$arr[0,19] = 0.4531346620072664
# row 3: ECs -> FAPs
$arr[1,0] = 'ECs'
$arr[1,1] = 'Sema3f'
$arr[1,2] = 'Nrp2'
$arr[1,3] = 'FAPs'
$arr[1,4] = 3
$arr[1,5] = 1
$arr[1,6] = 37.05582933333334
$arr[1,7] = 111.167488
$arr[1,8] = 0.9062076988734117
$arr[1,9] = 0.9062076988734117
$arr[1,10] = 3
$arr[1,11] = 1
$arr[1,12] = 8.622273333333332
$arr[1,13] = 25.86682
$arr[1,14] = 0.1362955132808722
$arr[1,15] = 0.1362955132808722
$arr[1,16] = 319.5054891053511
$arr[1,17] = 2875.54940194816
$arr[1,18] = 0.1235120434570297
$arr[1,19] = 0.1235120434570297
# row 4: ECs -> M2
$arr[2,0] = 'ECs'
$arr[2,1] = 'Sema3f'
$arr[2,2] = 'Nrp2'
$arr[2,3] = 'M2'
$arr[2,4] = 3
$arr[2,5] = 1
$arr[2,6] = 37.05582933333334
$arr[2,7] = 111.167488
$arr[2,8] = 0.9062076988734117
$arr[2,9] = 0.9062076988734117
$arr[2,10] = 3
$arr[2,11] = 1
$arr[2,12] = 18.70876033333333
$arr[2,13] = 56.12628100000001
$arr[2,14] = 0.2957364019791172
$arr[2,15] = 0.2957364019791172
$arr[2,16] = 693.2686299502365
$arr[2,17] = 6239.417669552129
$arr[2,18] = 0.267998604310598
$arr[2,19] = 0.2679986043105981
# row 5: ECs -> sCs
$arr[3,0] = 'ECs'
$arr[3,1] = 'Sema3f'
$arr[3,2] = 'Nrp2'
$arr[3,3] = 'sCs'
$arr[3,4] = 3
$arr[3,5] = 1
$arr[3,6] = 37.05582933333334
$arr[3,7] = 111.167488
$arr[3,8] = 0.9062076988734117
$arr[3,9] = 0.9062076988734117
$arr[3,10] = 3
$arr[3,11] = 1
$arr[3,12] = 4.297619333333333
$arr[3,13] = 12.892858
$arr[3,14] = 0.06793408307505136
$arr[3,15] = 0.06793408307505136
$arr[3,16] = 159.2518485556338
$arr[3,17] = 1433.266637000704
$arr[3,18] = 0.06156238909851747
$arr[3,19] = 0.06156238909851747
# row 6: FAPs -> ECs
$arr[4,0] = 'FAPs'
$arr[4,1] = 'Sema3f'
$arr[4,2] = 'Nrp2'
$arr[4,3] = 'ECs'
$arr[4,4] = 3
$arr[4,5] = 1
$arr[4,6] = 1.679068333333333
$arr[4,7] = 5.037205
$arr[4,8] = 0.04106195106076017
$arr[4,9] = 0.04106195106076017
$arr[4,10] = 3
$arr[4,11] = 1
$arr[4,12] = 31.632955
$arr[4,13] = 94.898865
$arr[4,14] = 0.5000340016649593
$arr[4,15] = 0.5000340016649593
$arr[4,16] = 53.11389303025833
$arr[4,17] = 478.025037272325
$arr[4,18] = 0.02053237170508263
$arr[4,19] = 0.02053237170508263
# row 7: FAPs -> FAPs
$arr[5,0] = 'FAPs'
$arr[5,1] = 'Sema3f'
$arr[5,2] = 'Nrp2'
$arr[5,3] = 'FAPs'
$arr[5,4] = 3
$arr[5,5] = 1
$arr[5,6] = 1.679068333333333
$arr[5,7] = 5.037205
$arr[5,8] = 0.04106195106076017
$arr[5,9] = 0.04106195106076017
$arr[5,10] = 3
$arr[5,11] = 1
$arr[5,12] = 8.622273333333332
$arr[5,13] = 25.86682
$arr[5,14] = 0.1362955132808722
$arr[5,15] = 0.1362955132808722
$arr[5,16] = 14.47738611534444
$arr[5,17] = 130.2964750381
$arr[5,18] = 0.005596559696140363
$arr[5,19] = 0.005596559696140363
# row 8: FAPs -> M2
$arr[6,0] = 'FAPs'
$arr[6,1] = 'Sema3f'
$arr[6,2] = 'Nrp2'
$arr[6,3] = 'M2'
$arr[6,4] = 3
$arr[6,5] = 1
$arr[6,6] = 1.679068333333333
$arr[6,7] = 5.037205
$arr[6,8] = 0.04106195106076017
$arr[6,9] = 0.04106195106076017
$arr[6,10] = 3
$arr[6,11] = 1
$arr[6,12] = 18.70876033333333
$arr[6,13] = 56.12628100000001
$arr[6,14] = 0.2957364019791172
$arr[6,15] = 0.2957364019791172
$arr[6,16] = 31.41328703162278
$arr[6,17] = 282.7195832846051
$arr[6,18] = 0.01214351366495181
$arr[6,19] = 0.01214351366495181
# row 9: FAPs -> sCs
$arr[7,0] = 'FAPs'
$arr[7,1] = 'Sema3f'
$arr[7,2] = 'Nrp2'
$arr[7,3] = 'sCs'
$arr[7,4] = 3
$arr[7,5] = 1
$arr[7,6] = 1.679068333333333
$arr[7,7] = 5.037205
$arr[7,8] = 0.04106195106076017
$arr[7,9] = 0.04106195106076017
$arr[7,10] = 3
$arr[7,11] = 1
$arr[7,12] = 4.297619333333333
$arr[7,13] = 12.892858
$arr[7,14] = 0.06793408307505136
$arr[7,15] = 0.06793408307505136
$arr[7,16] = 7.215996531321111
$arr[7,17] = 64.94396878189
$arr[7,18] = 0.002789505994585375
$arr[7,19] = 0.002789505994585375
# row 10: M2 -> ECs
$arr[8,0] = 'M2'
$arr[8,1] = 'Sema3f'
$arr[8,2] = 'Nrp2'
$arr[8,3] = 'ECs'
$arr[8,4] = 2
$arr[8,5] = 0.6666666666666666
$arr[8,6] = 0.08951333333333333
$arr[8,7] = 0.26854
$arr[8,8] = 0.002189066424308031
$arr[8,9] = 0.002189066424308031
$arr[8,10] = 3
$arr[8,11] = 1
$arr[8,12] = 31.632955
$arr[8,13] = 94.898865
$arr[8,14] = 0.5000340016649593
$arr[8,15] = 0.5000340016649593
$arr[8,16] = 2.831571245233333
$arr[8,17] = 25.4841412071
$arr[8,18] = 0.001094607644057149
$arr[8,19] = 0.001094607644057148
# row 11: M2 -> FAPs
$arr[9,0] = 'M2'
$arr[9,1] = 'Sema3f'
$arr[9,2] = 'Nrp2'
$arr[9,3] = 'FAPs'
$arr[9,4] = 2
$arr[9,5] = 0.6666666666666666
$arr[9,6] = 0.08951333333333333
$arr[9,7] = 0.26854
$arr[9,8] = 0.002189066424308031
$arr[9,9] = 0.002189066424308031
$arr[9,10] = 3
$arr[9,11] = 1
$arr[9,12] = 8.622273333333332
$arr[9,13] = 25.86682
$arr[9,14] = 0.1362955132808722
$arr[9,15] = 0.1362955132808722
$arr[9,16] = 0.7718084269777777
$arr[9,17] = 6.9462758428
$arr[9,18] = 0.0002983599319069867
$arr[9,19] = 0.0002983599319069867
# row 12: M2 -> M2
$arr[10,0] = 'M2'
$arr[10,1] = 'Sema3f'
$arr[10,2] = 'Nrp2'
$arr[10,3] = 'M2'
$arr[10,4] = 2
$arr[10,5] = 0.6666666666666666
$arr[10,6] = 0.08951333333333333
$arr[10,7] = 0.26854
$arr[10,8] = 0.002189066424308031
$arr[10,9] = 0.002189066424308031
$arr[10,10] = 3
$arr[10,11] = 1
$arr[10,12] = 18.70876033333333
$arr[10,13] = 56.12628100000001
$arr[10,14] = 0.2957364019791172
$arr[10,15] = 0.2957364019791172
$arr[10,16] = 1.674683499971111
$arr[10,17] = 15.07215149974
$arr[10,18] = 0.0006473866280181486
$arr[10,19] = 0.0006473866280181486
# row 13: M2 -> sCs
$arr[11,0] = 'M2'
$arr[11,1] = 'Sema3f'
$arr[11,2] = 'Nrp2'
$arr[11,3] = 'sCs'
$arr[11,4] = 2
$arr[11,5] = 0.6666666666666666
$arr[11,6] = 0.08951333333333333
$arr[11,7] = 0.26854
$arr[11,8] = 0.002189066424308031
$arr[11,9] = 0.002189066424308031
$arr[11,10] = 3
$arr[11,11] = 1
$arr[11,12] = 4.297619333333333
$arr[11,13] = 12.892858
$arr[11,14] = 0.06793408307505136
$arr[11,15] = 0.06793408307505136
$arr[11,16] = 0.3846942319244445
$arr[11,17] = 3.46224808732
$arr[11,18] = 0.0001487122203257474
$arr[11,19] = 0.0001487122203257474
# row 14: sCs -> ECs
$arr[12,0] = 'sCs'
$arr[12,1] = 'Sema3f'
$arr[12,2] = 'Nrp2'
$arr[12,3] = 'ECs'
$arr[12,4] = 3
$arr[12,5] = 1
$arr[12,6] = 2.066688666666666
$arr[12,7] = 6.200066
$arr[12,8] = 0.05054128364152006
$arr[12,9] = 0.05054128364152006
$arr[12,10] = 3
$arr[12,11] = 1
$arr[12,12] = 31.632955
$arr[12,13] = 94.898865
$arr[12,14] = 0.5000340016649593
$arr[12,15] = 0.5000340016649593
$arr[12,16] = 65.37546959167666
$arr[12,17] = 588.3792263250899
$arr[12,18] = 0.02527236030855302
$arr[12,19] = 0.02527236030855302
# row 15: sCs -> FAPs
$arr[13,0] = 'sCs'
$arr[13,1] = 'Sema3f'
$arr[13,2] = 'Nrp2'
$arr[13,3] = 'FAPs'
$arr[13,4] = 3
$arr[13,5] = 1
$arr[13,6] = 2.066688666666666
$arr[13,7] = 6.200066
$arr[13,8] = 0.05054128364152006
$arr[13,9] = 0.05054128364152006
$arr[13,10] = 3
$arr[13,11] = 1
$arr[13,12] = 8.622273333333332
$arr[13,13] = 25.86682
$arr[13,14] = 0.1362955132808722
$arr[13,15] = 0.1362955132808722
$arr[13,16] = 17.81955457890222
$arr[13,17] = 160.37599121012
$arr[13,18] = 0.006888550195795127
$arr[13,19] = 0.006888550195795127
# row 16: sCs -> M2
$arr[14,0] = 'sCs'
$arr[14,1] = 'Sema3f'
$arr[14,2] = 'Nrp2'
$arr[14,3] = 'M2'
$arr[14,4] = 3
$arr[14,5] = 1
$arr[14,6] = 2.066688666666666
$arr[14,7] = 6.200066
$arr[14,8] = 0.05054128364152006
$arr[14,9] = 0.05054128364152006
$arr[14,10] = 3
$arr[14,11] = 1
$arr[14,12] = 18.70876033333333
$arr[14,13] = 56.12628100000001
$arr[14,14] = 0.2957364019791172
$arr[14,15] = 0.2957364019791172
$arr[14,16] = 38.66518294828288
$arr[14,17] = 347.986646534546
$arr[14,18] = 0.01494689737554915
$arr[14,19] = 0.01494689737554916
# row 17: sCs -> sCs
$arr[15,0] = 'sCs'
$arr[15,1] = 'Sema3f'
$arr[15,2] = 'Nrp2'
$arr[15,3] = 'sCs'
$arr[15,4] = 3
$arr[15,5] = 1
$arr[15,6] = 2.066688666666666
$arr[15,7] = 6.200066
$arr[15,8] = 0.05054128364152006
$arr[15,9] = 0.05054128364152006
$arr[15,10] = 3
$arr[15,11] = 1
$arr[15,12] = 4.297619333333333
$arr[15,13] = 12.892858
$arr[15,14] = 0.06793408307505136
$arr[15,15] = 0.06793408307505136
$arr[15,16] = 8.881841169847554
$arr[15,17] = 79.936570528628
$arr[15,18] = 0.003433475761622758
$arr[15,19] = 0.003433475761622758

$ws.Range("A2:T17").Value = $arr
